$d = $word.ActiveDocument

# Locate the "%%SurveyName%%" placeholder and grab its paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("%%SurveyName%%")
if (-not $found) {
    throw "Could not find '%%SurveyName%%' placeholder"
}
$para = $rng.Paragraphs(1).Range

# Rebuild the paragraph so the "%%...%%" placeholder is split into three
# runs - the literal "%%" delimiters plus the inner token - matching the
# pattern already used elsewhere in this document (e.g. "%%Staff%%",
# "%%SLK%%") and renaming the token from SurveyName to AssessmentType.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="66587A54" w14:textId="3F5206ED" w:rsidR="002A3BA9" w:rsidRPr="00FE2F29" w:rsidRDefault="00FE2F29" w:rsidP="00FE2F29">
            <w:pPr>
              <w:pStyle w:val="NoSpacing"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
              </w:rPr>
              <w:t>%%</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
              </w:rPr>
              <w:t>AssessmentType</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
              </w:rPr>
              <w:t>%%</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$para.InsertXML($xml)
